$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$core = $wb.Worksheets.Item("Core")
$chart = $wb.Worksheets.Item("Test Chart")

# --- Style swap -------------------------------------------------------
# Originally only Metadata!G7 uses the "alignment" style; move that same
# style onto "Test Chart"!A2, and reset Metadata!G7 back to the plain
# style (matching Metadata!F7), by copying formats between cells so the
# underlying style table / cell values stay untouched.
$meta.Range("G7").Copy()
$chart.Range("A2").PasteSpecial(-4122)

$meta.Range("F7").Copy()
$meta.Range("G7").PasteSpecial(-4122)

# --- Core sheet: collapse placeholder "code" column onto the question's
# real code (column B) for the four core questions ---------------------
$core.Range("A2").Value = $core.Range("B2").Value2
$core.Range("A3").Value = $core.Range("B3").Value2
$core.Range("A4").Value = $core.Range("B4").Value2
$core.Range("A5").Value = $core.Range("B5").Value2

# --- Test Chart sheet: rename first question code ----------------------
$chart.Range("A2").Value = "PatientChartingDate"
